# "global switch parameters change from drop duplicates to Princiapl switch"
# 1) Sheet1 ("parameters"): a new customer block (3 columns: Complight, Gs-labs,
#    Rossiya) is inserted before the existing column C, pushing the previous
#    per-customer columns (Rencredit, IBS, Novatrans, ...) three columns to the
#    right.
# 2) Sheet2 ("report_info"): two "enabled" flags (export_to_excel /
#    force_extract) are flipped so that the "Principal switch" aggregated
#    report (isl_aggregated, row 50) is turned off and the "global fabric
#    parameters" report (row 47) is turned on instead; row 12 flag also
#    cleared.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: insert 3 new columns at C, fill in the new "Complight" / "Gs-labs" /
# "Rossiya" customer parameter block.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Columns("C:E").Insert()

# Column widths for the freshly inserted columns (C keeps the A:B width, D/E
# get their own custom widths).
$ws1.Columns("C").ColumnWidth = 29.7109375
$ws1.Columns("D").ColumnWidth = 52.140625
$ws1.Columns("E").ColumnWidth = 29.7109375

# The small "name/value" legend in row 1 only occupied A1/D1 before; after the
# insert the shifted-out copy (now at G1) is removed and a fresh "value"
# label is written at C1 (same fill/bold formatting as A1).
$ws1.Range("G1").Clear()
$ws1.Range("C1").Font.Bold = $true
$ws1.Range("C1").Interior.Color = 255
$ws1.Range("C1").Value = "value"

# New "Rossiya" column (E)
$ws1.Range("E2").Value = "Rossiya"
$ws1.Range("E4").Value = "C:\Users\vlasenko\Documents\06.CONFIGS\Rossiya\SANSW"
$ws1.Range("E3").Value = "C:\Users\vlasenko\Documents\01.CUSTOMERS\Rossiya"

# New "Gs-labs" column (D)
$ws1.Range("D4").Value = "C:\Users\vlasenko\Documents\06.CONFIGS\GS_labs\0212\92"
$ws1.Range("D2").Value = "Gs-labs"
$ws1.Range("D3").Value = "C:\Users\vlasenko\Documents\01.CUSTOMERS\Gs-labs"

# New "Complight" column (C)
$ws1.Range("C2").Value = "Complight"
$ws1.Range("C4").Value = "C:\Users\vlasenko\Documents\06.CONFIGS\Complight\ssave"
$ws1.Range("C3").Value = "C:\Users\vlasenko\Documents\01.CUSTOMERS\Complight\DEC2020"

[void]$ws1.Range("C16").Select()

# ---------------------------------------------------------------------------
# Sheet2: flip the export_to_excel / force_extract switches.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("G12").Value = 0

$ws2.Range("G47").Value = 1
$ws2.Range("H47").Value = 1

$ws2.Range("G50").Value = 0
$ws2.Range("H50").Value = 0

[void]$ws2.Activate()
[void]$ws2.Range("G52").Select()
